$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'255.82"
$ws.Range("E2").Value = "'4.15%"
$ws.Range("D3").Value = "'28.09"
$ws.Range("E3").Value = "'-4.48%"
$ws.Range("D4").Value = "'5.352"
$ws.Range("E4").Value = "'3.95%"
$ws.Range("D5").Value = "'0.05817"
$ws.Range("E5").Value = "'0.70%"
$ws.Range("D6").Value = "'6.707"
$ws.Range("E6").Value = "'1.36%"
$ws.Range("D7").Value = "'3.248"
$ws.Range("E7").Value = "'2.61%"
$ws.Range("D8").Value = "'0.8711"
$ws.Range("E8").Value = "'1.69%"
$ws.Range("D9").Value = "'0.8957"
$ws.Range("E9").Value = "'4.55%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1416"
$ws.Range("E10").Value = "'3.94%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07242"
$ws.Range("E11").Value = "'3.09%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03179"
$ws.Range("E12").Value = "'5.01%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09247"
$ws.Range("E13").Value = "'-1.33%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001540"
$ws.Range("E14").Value = "'1.04%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006057"
$ws.Range("E15").Value = "'-94.07%"
$ws.Range("D16").Value = "'0.006031"
$ws.Range("E16").Value = "'1.01%"
$ws.Range("D17").Value = "'3.501"
$ws.Range("E17").Value = "'0.46%"
$ws.Range("E18").Value = "'4.70%"
$ws.Range("E19").Value = "'-1.14%"
$ws.Range("D20").Value = "'0.03440"
$ws.Range("E20").Value = "'3.58%"
$ws.Range("E21").Value = "'2.41%"
$ws.Range("D22").Value = "'3.522"
$ws.Range("E22").Value = "'6.13%"
$ws.Range("D23").Value = "'0.04163"
$ws.Range("E23").Value = "'0.51%"
$ws.Range("E24").Value = "'-1.49%"
$ws.Range("D25").Value = "'0.001222"
$ws.Range("E25").Value = "'-0.29%"
$ws.Range("D26").Value = "'0.004875"
$ws.Range("E26").Value = "'17.92%"
$ws.Range("E27").Value = "'-0.81%"
$ws.Range("E28").Value = "'0.70%"
$ws.Range("D41").Value = "'0.005756"
$ws.Range("E41").Value = "'-2.05%"
$ws.Range("E42").Value = "'3.01%"
$ws.Range("D43").Value = "'0.002199"
$ws.Range("E43").Value = "'0.00%"
$ws.Range("D44").Value = "'0.009927"
$ws.Range("E44").Value = "'16.43%"
$ws.Range("D45").Value = "'0.00005284"
$ws.Range("E45").Value = "'0.15%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("D47").Value = "'0.08496"
$ws.Range("E47").Value = "'46.56%"
$ws.Range("D48").Value = "'0.002140"
$ws.Range("E48").Value = "'-1.42%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'0.00%"
